$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 565.4286
$ws.Range("I2").Value = 695
$ws.Range("J2").Value = 392.66666
$ws.Range("K2").Value = 695
$ws.Range("L2").Value = 392.66666
$ws.Range("M2").Value = -582
$ws.Range("N2").Value = -618.66666
$ws.Range("H106").Value = 1805.8889
$ws.Range("I106").Value = 875.5
$ws.Range("J106").Value = 3666.6667
$ws.Range("K106").Value = 875.5
$ws.Range("L106").Value = 3666.6667
$ws.Range("M106").Value = -244.5
$ws.Range("N106").Value = -4928.6667
$ws.Range("H107").Value = 1070.9286
$ws.Range("I107").Value = 1123.9166
$ws.Range("K107").Value = 1123.9166
$ws.Range("M107").Value = 796.0834
$ws.Range("H113").Value = 8024.75
$ws.Range("I113").Value = 3654.4443
$ws.Range("J113").Value = 13643.714
$ws.Range("K113").Value = 3654.4443
$ws.Range("L113").Value = 13643.714
$ws.Range("M113").Value = -400.4443000000001
$ws.Range("N113").Value = -20151.714
$ws.Range("H116").Value = 6850.3
$ws.Range("I116").Value = 1808.3334
$ws.Range("K116").Value = 1808.3334
$ws.Range("M116").Value = 1633.6666
$ws.Range("H132").Value = 89781.78
$ws.Range("I132").Value = 107077.945
$ws.Range("J132").Value = 7625
$ws.Range("K132").Value = 321233.835
$ws.Range("L132").Value = 22875
$ws.Range("M132").Value = -318703.835
$ws.Range("N132").Value = -27935
$ws.Range("H137").Value = 4585.3716
$ws.Range("I137").Value = 4029.4
$ws.Range("K137").Value = 12088.2
$ws.Range("M137").Value = -9538.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 753.1739
$ws.Range("I2").Value = 774.3333
$ws.Range("J2").Value = 677
$ws.Range("K2").Value = 774.3333
$ws.Range("L2").Value = 677
$ws.Range("M2").Value = -661.3333
$ws.Range("N2").Value = -903
$ws.Range("H32").Value = 7999.0933
$ws.Range("I32").Value = 5073.7
$ws.Range("J32").Value = 11342.4
$ws.Range("K32").Value = 5073.7
$ws.Range("L32").Value = 11342.4
$ws.Range("M32").Value = -4786.7
$ws.Range("N32").Value = -11916.4
$ws.Range("H61").Value = 2098.0417
$ws.Range("I61").Value = 1135.0769
$ws.Range("K61").Value = 1135.0769
$ws.Range("M61").Value = -923.0769
$ws.Range("H116").Value = 753.1739
$ws.Range("I116").Value = 774.3333
$ws.Range("J116").Value = 677
$ws.Range("K116").Value = 774.3333
$ws.Range("L116").Value = 677
$ws.Range("M116").Value = 1519.6667
$ws.Range("N116").Value = -5265
$ws.Range("H136").Value = 2098.0417
$ws.Range("I136").Value = 1135.0769
$ws.Range("K136").Value = 3405.2307
$ws.Range("M136").Value = -855.2307000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 753.1739
$ws.Range("I3").Value = 774.3333
$ws.Range("J3").Value = 677
$ws.Range("K3").Value = 774.3333
$ws.Range("L3").Value = 677
$ws.Range("M3").Value = -660.3333
$ws.Range("N3").Value = -905
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H94").Value = 1490.2632
$ws.Range("I94").Value = 1399.375
$ws.Range("J94").Value = 1975
$ws.Range("K94").Value = 1399.375
$ws.Range("L94").Value = 1975
$ws.Range("M94").Value = -948.375
$ws.Range("N94").Value = -2877
$ws.Range("H138").Value = 41295.715
$ws.Range("J138").Value = 41295.715
$ws.Range("L138").Value = 41295.715
$ws.Range("N138").Value = -51575.715
$ws.Range("H140").Value = 48168.57
$ws.Range("J140").Value = 48168.57
$ws.Range("L140").Value = 48168.57
$ws.Range("N140").Value = -58528.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3641.8857
$ws.Range("I31").Value = 1908.8889
$ws.Range("J31").Value = 4241.769
$ws.Range("K31").Value = 1908.8889
$ws.Range("L31").Value = 4241.769
$ws.Range("M31").Value = -1613.8889
$ws.Range("N31").Value = -4831.769
$ws.Range("H34").Value = 3641.8857
$ws.Range("I34").Value = 1908.8889
$ws.Range("J34").Value = 4241.769
$ws.Range("K34").Value = 1908.8889
$ws.Range("L34").Value = 4241.769
$ws.Range("M34").Value = -1706.8889
$ws.Range("N34").Value = -4645.769
$ws.Range("H68").Value = 53282.125
$ws.Range("J68").Value = 53282.125
$ws.Range("L68").Value = 53282.125
$ws.Range("N68").Value = -54780.125
$ws.Range("H71").Value = 53282.125
$ws.Range("J71").Value = 53282.125
$ws.Range("L71").Value = 159846.375
$ws.Range("N71").Value = -167334.375
$ws.Range("H132").Value = 4022.4517
$ws.Range("I132").Value = 3317.5881
$ws.Range("J132").Value = 4878.357
$ws.Range("K132").Value = 9952.764299999999
$ws.Range("L132").Value = 14635.071
$ws.Range("M132").Value = -7422.764299999999
$ws.Range("N132").Value = -19695.071
$ws.Range("H134").Value = 7068.952
$ws.Range("I134").Value = 8338.714
$ws.Range("J134").Value = 4529.4287
$ws.Range("K134").Value = 25016.142
$ws.Range("L134").Value = 13588.2861
$ws.Range("M134").Value = -22481.142
$ws.Range("N134").Value = -18658.2861
$ws.Range("H141").Value = 15722.728
$ws.Range("J141").Value = 15722.728
$ws.Range("L141").Value = 15722.728
$ws.Range("N141").Value = -26082.728

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5942.7856
$ws.Range("I56").Value = 5942.7856
$ws.Range("K56").Value = 5942.7856
$ws.Range("M56").Value = -5412.7856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6590.968
$ws.Range("I70").Value = 5855
$ws.Range("K70").Value = 5855
$ws.Range("M70").Value = -5585
$ws.Range("H73").Value = 6590.968
$ws.Range("I73").Value = 5855
$ws.Range("K73").Value = 5855
$ws.Range("M73").Value = -4919
$ws.Range("H126").Value = 4185.3657
$ws.Range("I126").Value = 2970.7317
$ws.Range("J126").Value = 5400
$ws.Range("K126").Value = 8912.195099999999
$ws.Range("L126").Value = 16200
$ws.Range("M126").Value = -6442.195099999999
$ws.Range("N126").Value = -21140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1738
$ws.Range("I100").Value = 1618.8
$ws.Range("J100").Value = 2334
$ws.Range("K100").Value = 1618.8
$ws.Range("L100").Value = 2334
$ws.Range("M100").Value = -1077.8
$ws.Range("N100").Value = -3416
$ws.Range("H136").Value = 4801.846
$ws.Range("I136").Value = 2267.7144
$ws.Range("J136").Value = 7758.3335
$ws.Range("K136").Value = 6803.1432
$ws.Range("L136").Value = 23275.0005
$ws.Range("M136").Value = -4253.1432
$ws.Range("N136").Value = -28375.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
